$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(108, 2017, 99999, 1, "[4942]", 4942, 5988, 0.8253),
    @(109, 2022, 349999, 1, "[4934]", 4934, 5975, 0.8258),
    @(110, 2026, 49999, 1, "[4938]", 4938, 6006, 0.8222),
    @(111, 2030, 74999, 1, "[4960]", 4960, 6043, 0.8208),
    @(112, 2034, 149999, 1, "[5024]", 5024, 6000, 0.8373),
    @(113, 2038, 99999, 1, "[4992]", 4992, 5978, 0.8351),
    @(114, 2042, 349999, 1, "[5039]", 5039, 6052, 0.8326),
    @(115, 2046, 149999, 1, "[5061]", 5061, 5966, 0.8483000000000001),
    @(116, 2050, 274999, 1, "[4963]", 4963, 5986, 0.8290999999999999),
    @(117, 2054, 99999, 1, "[4957]", 4957, 5953, 0.8327),
    @(118, 2058, 199999, 1, "[4962]", 4962, 5950, 0.8339),
    @(119, 2062, 149999, 1, "[4939]", 4939, 5948, 0.8304),
    @(120, 2066, 74999, 1, "[4961]", 4961, 5996, 0.8274),
    @(121, 2070, 149999, 1, "[4992]", 4992, 5985, 0.8341),
    @(122, 2074, 399999, 1, "[4938]", 4938, 5887, 0.8388),
    @(123, 2078, 124999, 1, "[4999]", 4999, 6037, 0.8280999999999999),
    @(124, 2082, 49999, 1, "[4933]", 4933, 5937, 0.8309),
    @(125, 2086, 124999, 1, "[5058]", 5058, 6064, 0.8341),
    @(126, 2090, 124999, 1, "[4942]", 4942, 5907, 0.8366),
    @(127, 2094, 74999, 1, "[4949]", 4949, 5968, 0.8293),
    @(128, 2098, 124999, 1, "[4961]", 4961, 6012, 0.8252),
    @(129, 2102, 24999, 1, "[5049]", 5049, 6084, 0.8299),
    @(130, 2106, 49999, 1, "[5004]", 5004, 5997, 0.8344),
    @(131, 2110, 149999, 1, "[5033]", 5033, 5983, 0.8411999999999999),
    @(132, 2114, 199999, 1, "[4967]", 4967, 5995, 0.8285),
    @(133, 2118, 99999, 1, "[5038]", 5038, 6065, 0.8307),
    @(134, 2122, 224999, 1, "[4980]", 4980, 6039, 0.8246),
    @(135, 2126, 174999, 1, "[4957]", 4957, 5923, 0.8369),
    @(136, 2130, 99999, 1, "[4948]", 4948, 6028, 0.8208),
    @(137, 2134, 99999, 1, "[4945]", 4945, 5986, 0.8260999999999999),
    @(138, 2138, 224999, 1, "[5005]", 5005, 6016, 0.8319)
)

foreach ($r in $data) {
    $rowNum = $r[0]
    $ws.Cells.Item($rowNum, 1).Value = $r[1]
    $ws.Cells.Item($rowNum, 2).Value = $r[2]
    $ws.Cells.Item($rowNum, 3).Value = $r[3]
    $ws.Cells.Item($rowNum, 4).Value = $r[4]
    $ws.Cells.Item($rowNum, 5).Value = $r[5]
    $ws.Cells.Item($rowNum, 6).Value = $r[6]
    $ws.Cells.Item($rowNum, 7).Value = $r[7]
}
